# Update the "Start Date" row (row 4) values for Senate (B), House (C), and
# Total (D) columns from 1/3/1993 to 1/5/1993 on the "103_1" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("103_1")

# 33974 is the Excel date serial number for 1/5/1993 (days since 12/30/1899),
# matching the date format (m/d/yyyy) already applied to these cells.
$newDateSerial = 33974

$ws.Range("B4").Value = $newDateSerial
$ws.Range("C4").Value = $newDateSerial
$ws.Range("D4").Value = $newDateSerial
